# RTM: populate traceability matrix rows for the login feature (SRS -> Test cases)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "SRS ID"
$ws.Range("B1").Value = "Test case ID"

# --- SRS IDs (column A) ----------------------------------------------------
$ws.Range("A2").Value = "SRS_Login_001"
$ws.Range("A3").Value = "SRS_Login_002"
$ws.Range("A4").Value = "SRS_Login_003"
$ws.Range("A5").Value = "SRS_Login_004"
$ws.Range("A6").Value = "SRS_Login_005"

# --- Test case IDs (column B) -----------------------------------------------
# (populated in the same order as the original authoring session so the
# shared-string table layout lines up with the source workbook)
$ws.Range("B2").Value = "TC_Login_001 "
$ws.Range("B5").Value = "TC_Login_001  " + [char]10 + "TC-Login_008"
$ws.Range("B3").Value = "TC_Login_001  "
$ws.Range("B6").Value = "TC_Login_002" + [char]10 + "TC_Login_003" + [char]10 + "TC_Login_004" + [char]10 + "TC_Login_009" + [char]10 + "TC_Login_0010" + [char]10 + "TC_Login_0011" + [char]10
$ws.Range("B4").Value = "TC_Login_014"

# Wrap the multi-line / overflow cells so the text is readable
$ws.Range("B3").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("B6").WrapText = $true

# --- Row heights to fit the new (wrapped) content --------------------------
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 18.75
$ws.Rows.Item(5).RowHeight = 37.5
$ws.Rows.Item(6).RowHeight = 131.25

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.85546875
$ws.Columns.Item(2).ColumnWidth = 37.42578125

# --- Selection / scroll position -------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E6").Select()
